$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on D-column price cells being updated, so numeric-looking
# strings like "1.00" or "211.39" stay text (matching the original "inlineStr" cells)
# instead of Excel auto-converting them to numbers on assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.694.85"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.599.83"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "211.39"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "19.56"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.824.21"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.05"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.565.09"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "65.36"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "26.673.95"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "0.0₃0761"
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").Value = "210.18"
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "7.19"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "143.14"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("D34").Value = "1.291.06"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "1.06"
$ws.Range("E39").Value = "  +15.86%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D44").Value = "63.21"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "1.735.05"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "91.24"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "7.36"
$ws.Range("E51").Value = "  -1.07%  "
